$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.176679611206055
$ws.Range("B1").Value = 2.395702600479126
$ws.Range("C1").Value = 6.489742279052734
$ws.Range("D1").Value = 2.061849594116211
$ws.Range("E1").Value = 1.199370980262756
